# Swap the "Enterprises (absolute #)" row (row 12) with the
# "Enterprises density (per 1000 people)" row (row 13) on the Summary
# sheet, so the density figures now appear before the absolute-count
# figures (matching the reordering of the shared-strings table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Capture current (pre-edit) values.
# NOTE: in this COM-interop runtime, ".Value" (as a plain property read)
# yields a parameterized-property descriptor rather than the cell's
# value; it must be invoked as a method, ".Value()", to get the actual
# contents back.
$A12 = $ws.Range("A12").Value()
$C12 = $ws.Range("C12").Value()
$D12 = $ws.Range("D12").Value()

$A13 = $ws.Range("A13").Value()
$C13 = $ws.Range("C13").Value()
$D13 = $ws.Range("D13").Value()

# The numeric-looking labels ("39.8" / "2646549") are stored as *text*
# in the workbook, not numbers. Force the destination cells to Text
# format before assigning so the runtime keeps them as text instead of
# silently converting them to numeric values.
$ws.Range("C12").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"

# Row 12 now gets what used to be row 13's content (the density row).
$ws.Range("A12").Value = $A13
$ws.Range("C12").Value = $C13
$ws.Range("D12").Value = $D13

# Row 13 now gets what used to be row 12's content (the absolute-# row).
$ws.Range("A13").Value = $A12
$ws.Range("C13").Value = $C12
$ws.Range("D13").Value = $D12
